{"js": "// Load all paragraphs in the document body so we can find the two\n// TODO list items that need to be merged.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst oldFirstText = \"Add code to move eyes back to square in front of box(home).\";\nconst secondText = \"Add animation to move the eyes in the box.\";\n\nlet firstPara = null;\nlet secondPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = (p.text || \"\").trim();\n  if (firstPara === null && t === oldFirstText) {\n    firstPara = p;\n    continue;\n  }\n  if (secondPara === null && t === secondText) {\n    secondPara = p;\n  }\n}\n\nif (firstPara && secondPara) {\n  // Replace the first paragraph's text with the second paragraph's text,\n  // then remove the now-duplicate second paragraph entirely.\n  firstPara.insertText(secondText, \"Replace\");\n  secondPara.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$oldFirstText = \"Add code to move eyes back to square in front of box(home).\"\n$secondText = \"Add animation to move the eyes in the box.\"\n\n# Replace the first TODO bullet's text with the second bullet's text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldFirstText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $secondText\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Now remove the now-duplicate bullet paragraph (the original second\n# paragraph, which already had this text before the replace above).\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $secondText) {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
